$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.121.47'
$ws.Range("E2").Value = '  +0.32%  '

$ws.Range("D3").Value = '2.328.88'
$ws.Range("E3").Value = '  +1.10%  '

$ws.Range("D5").Value = '''304.34'
$ws.Range("E5").Value = '  +1.51%  '

$ws.Range("D6").Value = '''97.89'
$ws.Range("E6").Value = '  +0.45%  '

$ws.Range("E7").Value = '  -1.29%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("D10").Value = '''35.63'
$ws.Range("E10").Value = '  -0.19%  '

$ws.Range("D11").Value = '''19.47'
$ws.Range("E11").Value = '  +7.93%  '

$ws.Range("E12").Value = '  +1.57%  '

$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("E14").Value = '  +1.75%  '

$ws.Range("D15").Value = '2.692.05'
$ws.Range("E15").Value = '  +1.27%  '

$ws.Range("D16").Value = '2.314.56'
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").Value = '43.029.75'
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("D19").Value = '''12.61'
$ws.Range("E19").Value = '  -0.83%  '

$ws.Range("D20").Value = '0.0₃0903'
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("D22").Value = '''68.03'
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").Value = '''237.50'
$ws.Range("E23").Value = '  -1.14%  '

$ws.Range("E24").Value = '  +3.98%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '''2.44'
$ws.Range("E26").Value = '  +0.60%  '

$ws.Range("D27").Value = '''24.96'
$ws.Range("E27").Value = '  -2.07%  '

$ws.Range("D28").Value = '''166.21'
$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("E29").Value = '  +2.14%  '

$ws.Range("E30").Value = '  +1.00%  '

$ws.Range("D31").Value = '''33.18'
$ws.Range("E31").Value = '  +0.15%  '

$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("D33").Value = '''17.99'
$ws.Range("E33").Value = '  +6.01%  '

$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("D35").Value = '''4.54'
$ws.Range("E35").Value = '  -8.44%  '

$ws.Range("E36").Value = '  +1.35%  '

$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("E38").Value = '  +0.30%  '

$ws.Range("D39").Value = '''2.80'
$ws.Range("E39").Value = '  +2.36%  '

$ws.Range("D40").Value = '''1.77'
$ws.Range("E40").Value = '  +0.18%  '

$ws.Range("D41").Value = '''0.110'
$ws.Range("E41").Value = '  -0.29%  '

$ws.Range("D42").Value = '1.997.39'
$ws.Range("E42").Value = '  -0.65%  '

$ws.Range("D43").Value = '''10.75'
$ws.Range("E43").Value = '  +6.03%  '

$ws.Range("E44").Value = '  -0.18%  '

$ws.Range("D45").Value = '''18.20'
$ws.Range("E45").Value = '  +5.42%  '

$ws.Range("E46").Value = '  -2.12%  '

$ws.Range("E47").Value = '  -0.33%  '

$ws.Range("D48").Value = '2.558.01'
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("D49").Value = '''2.86'
$ws.Range("E49").Value = '  -0.29%  '

$ws.Range("D50").Value = '''53.77'
$ws.Range("E50").Value = '  +0.33%  '

$ws.Range("D51").Value = '''71.98'
$ws.Range("E51").Value = '  -0.16%  '
